$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 428.78378
$ws.Range("J17").Value = 428.78378
$ws.Range("L17").Value = 1286.35134
$ws.Range("N17").Value = -1622.35134
$ws.Range("H40").Value = 3036.125
$ws.Range("J40").Value = 2069.8572
$ws.Range("L40").Value = 2069.8572
$ws.Range("N40").Value = -2419.8572
$ws.Range("H74").Value = 3666.6667
$ws.Range("I74").Value = 3000
$ws.Range("K74").Value = 3000
$ws.Range("M74").Value = -2064
$ws.Range("H77").Value = 3666.6667
$ws.Range("I77").Value = 3000
$ws.Range("K77").Value = 15000
$ws.Range("M77").Value = -10320
$ws.Range("H92").Value = 497.90475
$ws.Range("I92").Value = 368.73685
$ws.Range("K92").Value = 368.73685
$ws.Range("M92").Value = 879.26315
$ws.Range("H94").Value = 5500
$ws.Range("I94").Value = 5500
$ws.Range("K94").Value = 5500
$ws.Range("M94").Value = -5049
$ws.Range("H96").Value = 347.58823
$ws.Range("I96").Value = 238.09091
$ws.Range("J96").Value = 548.3333
$ws.Range("K96").Value = 714.27273
$ws.Range("L96").Value = 1644.9999
$ws.Range("M96").Value = 658.72727
$ws.Range("N96").Value = -4390.9999
$ws.Range("H112").Value = 2226.7026
$ws.Range("J112").Value = 2559.6
$ws.Range("L112").Value = 7678.799999999999
$ws.Range("N112").Value = -9894.799999999999
$ws.Range("H132").Value = 8776965
$ws.Range("I132").Value = 10758311
$ws.Range("K132").Value = 32274933
$ws.Range("M132").Value = -32272403
$ws.Range("H137").Value = 1241.069
$ws.Range("I137").Value = 924.15
$ws.Range("K137").Value = 2772.45
$ws.Range("M137").Value = -222.4499999999998
$ws.Range("H140").Value = 70000
$ws.Range("J140").Value = 70000
$ws.Range("L140").Value = 70000
$ws.Range("N140").Value = -80360
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3254.2097
$ws.Range("I32").Value = 3096.0908
$ws.Range("K32").Value = 3096.0908
$ws.Range("M32").Value = -2809.0908
$ws.Range("I63").Value = 1499.5
$ws.Range("J63").Value = 2100
$ws.Range("K63").Value = 1499.5
$ws.Range("L63").Value = 2100
$ws.Range("M63").Value = -813.5
$ws.Range("N63").Value = -3472
$ws.Range("I66").Value = 1499.5
$ws.Range("J66").Value = 2100
$ws.Range("K66").Value = 7497.5
$ws.Range("L66").Value = 10500
$ws.Range("M66").Value = -4065.5
$ws.Range("N66").Value = -17364
$ws.Range("H76").Value = 29998.5
$ws.Range("J76").Value = 29998.5
$ws.Range("L76").Value = 29998.5
$ws.Range("N76").Value = -30674.5
$ws.Range("H79").Value = 29998.5
$ws.Range("J79").Value = 29998.5
$ws.Range("L79").Value = 29998.5
$ws.Range("N79").Value = -32338.5
$ws.Range("H102").Value = 8774600
$ws.Range("I102").Value = 9806529
$ws.Range("J102").Value = 3200
$ws.Range("K102").Value = 9806529
$ws.Range("L102").Value = 3200
$ws.Range("M102").Value = -9804907
$ws.Range("N102").Value = -6444
$ws.Range("H132").Value = 2185.7878
$ws.Range("I132").Value = 1904.76
$ws.Range("J132").Value = 3064
$ws.Range("K132").Value = 5714.28
$ws.Range("L132").Value = 9192
$ws.Range("M132").Value = -3184.28
$ws.Range("N132").Value = -14252
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 76925600
$ws.Range("I105").Value = 90911740
$ws.Range("J105").Value = 1800
$ws.Range("K105").Value = 90911740
$ws.Range("L105").Value = 1800
$ws.Range("M105").Value = -90909993
$ws.Range("N105").Value = -5294
$ws.Range("H134").Value = 7827.8887
$ws.Range("I134").Value = 1301
$ws.Range("J134").Value = 24797.8
$ws.Range("K134").Value = 3903
$ws.Range("L134").Value = 74393.39999999999
$ws.Range("M134").Value = -1368
$ws.Range("N134").Value = -79463.39999999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 738.9091
$ws.Range("I105").Value = 712.8
$ws.Range("K105").Value = 712.8
$ws.Range("M105").Value = 1034.2
$ws.Range("H132").Value = 6401.0835
$ws.Range("I132").Value = 8354.666999999999
$ws.Range("K132").Value = 25064.001
$ws.Range("M132").Value = -22534.001
$ws.Range("H134").Value = 1525.25
$ws.Range("I134").Value = 1420.5938
$ws.Range("K134").Value = 4261.7814
$ws.Range("M134").Value = -1726.7814
$ws.Range("H141").Value = 29622.5
$ws.Range("J141").Value = 29622.5
$ws.Range("L141").Value = 29622.5
$ws.Range("N141").Value = -39982.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 228.22223
$ws.Range("I13").Value = 132.375
$ws.Range("K13").Value = 397.125
$ws.Range("M13").Value = -229.125
$ws.Range("H38").Value = 53.666668
$ws.Range("I38").Value = 53.666668
$ws.Range("K38").Value = 161.000004
$ws.Range("M38").Value = 185.999996
$ws.Range("H86").Value = 705
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 705
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 2115
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -4487
$ws.Range("H89").Value = 705
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 705
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 6345
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -18201
$ws.Range("H92").Value = 667.6
$ws.Range("I92").Value = 681.3333
$ws.Range("K92").Value = 2043.9999
$ws.Range("M92").Value = -795.9999
$ws.Range("H131").Value = 18870836
$ws.Range("J131").Value = 3186.9167
$ws.Range("L131").Value = 9560.750100000001
$ws.Range("N131").Value = -19640.7501
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 50003500
$ws.Range("I70").Value = 41670316
$ws.Range("K70").Value = 41670316
$ws.Range("M70").Value = -41670046
$ws.Range("H73").Value = 50003500
$ws.Range("I73").Value = 41670316
$ws.Range("K73").Value = 41670316
$ws.Range("M73").Value = -41669380
$ws.Range("H80").Value = 3857.75
$ws.Range("I80").Value = 3077
$ws.Range("J80").Value = 6200
$ws.Range("K80").Value = 3077
$ws.Range("L80").Value = 6200
$ws.Range("M80").Value = -2079
$ws.Range("N80").Value = -8196
$ws.Range("H83").Value = 3857.75
$ws.Range("I83").Value = 3077
$ws.Range("J83").Value = 6200
$ws.Range("K83").Value = 15385
$ws.Range("L83").Value = 31000
$ws.Range("M83").Value = -10393
$ws.Range("N83").Value = -40984
$ws.Range("H132").Value = 2071.6453
$ws.Range("I132").Value = 1873.7727
$ws.Range("K132").Value = 5621.3181
$ws.Range("M132").Value = -3091.3181
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H53").Value = 9525.5
$ws.Range("J53").Value = 9525.5
$ws.Range("L53").Value = 9525.5
$ws.Range("N53").Value = -10561.5
$ws.Range("H68").Value = 1972
$ws.Range("I68").Value = 1726.5
$ws.Range("K68").Value = 1726.5
$ws.Range("M68").Value = -977.5
$ws.Range("H71").Value = 1972
$ws.Range("I71").Value = 1726.5
$ws.Range("K71").Value = 8632.5
$ws.Range("M71").Value = -4888.5
$ws.Range("H137").Value = 36395.332
$ws.Range("J137").Value = 36395.332
$ws.Range("L137").Value = 36395.332
$ws.Range("N137").Value = -46595.332
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1471.8959
$ws.Range("I132").Value = 1002.54285
$ws.Range("K132").Value = 3007.62855
$ws.Range("M132").Value = -477.6285500000004
